$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) from 45174 to 45175
# for data rows 2 through 15.
$ws.Range("C2:C15").Value = 45175
